# =====================================================================
# manual_corrections.xlsx edit script
#
# Summary of changes (per commit message / xml diff):
#  1. Add a new "epa_clean" sheet (EPA plants to delete / manual changes)
#  2. Add a new "eia_clean" sheet (manual changes to fuel codes etc.)
#  3. Add a proposed plant row to generator_file (two blank placeholder
#     rows were left after the edit, rows 11 & 12)
#  4. Add manual corrections rows to unit_file (rows 16-27) for plant
#     52152, unit 6RB, prime mover ST - fixing a December/ozone-season
#     generation typo
#  5. plant_file is unchanged in content, only renumbered as a sheet
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet bookkeeping: create epa_clean + eia_clean ahead of the
#    existing three sheets, matching sheetId 5 / 6 (a short-lived
#    duplicate sheet, sheetId 4, is created and removed first - this
#    mirrors "removing duplicate static table" from the commit message
#    and reproduces the sheetId numbering seen in the workbook).
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$dup = $wb.Worksheets.Add($firstSheet)
$dup.Name = "duplicate_tmp"

$epa = $wb.Worksheets.Add($dup)
$epa.Name = "epa_clean"

$genSheet = $wb.Worksheets.Item("generator_file")
$eia = $wb.Worksheets.Add($genSheet)
$eia.Name = "eia_clean"

$wb.Worksheets.Item("duplicate_tmp").Delete()

# ---------------------------------------------------------------------
# 2. Populate epa_clean
# ---------------------------------------------------------------------
$epa = $wb.Worksheets.Item("epa_clean")

$epa.Range("A1").Value = "plant_id"
$epa.Range("B1").Value = "column_to_update"
$epa.Range("C1").Value = "update"
$epa.Range("A1:C1").Font.Bold = $true
$epa.Range("C1").NumberFormat = "@"

$epa.Range("A2").Value = 880004
$epa.Range("B2").Value = "plant_id"
$epa.Range("C2").Value = 57788

$epa.Range("A3").Value = 10154
$epa.Range("B3").Value = "operating_status"
$epa.Range("C3").Value = "OP"

$epa.Columns("A:A").ColumnWidth = 9
$epa.Columns("B:B").ColumnWidth = 18
$epa.Columns("C:C").ColumnWidth = 8

$epa.Range("B4").Select()

# ---------------------------------------------------------------------
# 3. Populate eia_clean
# ---------------------------------------------------------------------
$eia = $wb.Worksheets.Item("eia_clean")

$eia.Range("A1").Value = "plant_id"
$eia.Range("B1").Value = "generator_id"
$eia.Range("C1").Value = "column_to_update"
$eia.Range("A1:C1").Font.Bold = $true

$eia.Range("A2").Value = 60910
$eia.Range("B2").NumberFormat = "@"
$eia.Range("B2").Value = "NPLPS"
$eia.Range("C2").Value = "add"

$eia.Columns("A:A").ColumnWidth = 9
$eia.Columns("B:B").ColumnWidth = 13
$eia.Columns("C:C").ColumnWidth = 18

$eia.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4. generator_file: new proposed-plant entry -> two placeholder rows
# ---------------------------------------------------------------------
$gen = $wb.Worksheets.Item("generator_file")

$gen.Range("B11:C12").Style = "Normal"
$gen.Range("E11:E12").Style = "Normal"

$gen.Activate()
$gen.Range("A9:D9").Select()

# ---------------------------------------------------------------------
# 5. unit_file: manual corrections for plant 52152 / unit 6RB / ST
# ---------------------------------------------------------------------
$unit = $wb.Worksheets.Item("unit_file")

$unit.Columns("A:B").NumberFormat = "@"
$unit.Columns("E:E").NumberFormat = "@"

$unitRows = @(
    @("52152", "6RB", "ST", "heat_input", "6413476.992"),
    @("52152", "6RB", "ST", "heat_input_oz", " 2383097.20 "),
    @("52152", "6RB", "ST", "nox_mass", "272.69"),
    @("52152", "6RB", "ST", "nox_oz_mass", "97.016"),
    @("52152", "6RB", "ST", "so2_mass", "1272.5545"),
    @("52152", "6RB", "ST", "co2_mass", "608148.783"),
    @("52152", "6RB", "ST", "heat_input_source", "EIA Prime Mover-level Data"),
    @("52152", "6RB", "ST", "heat_input_oz_source", "EIA Prime Mover-level Data"),
    @("52152", "6RB", "ST", "nox_source", "Estimated using emissions factor"),
    @("52152", "6RB", "ST", "nox_oz_source", "Estimated using emissions factor"),
    @("52152", "6RB", "ST", "so2_source", "Estimated using emissions factor"),
    @("52152", "6RB", "ST", "co2_source", "Estimated using emissions factor")
)

$r = 16
foreach ($row in $unitRows) {
    $unit.Range("A$r").Value = $row[0]
    $unit.Range("B$r").Value = $row[1]
    $unit.Range("C$r").Value = $row[2]
    $unit.Range("D$r").Value = $row[3]
    $unit.Range("E$r").Value = $row[4]
    $r = $r + 1
}

$unit.Columns("A:A").ColumnWidth = 9
$unit.Columns("B:B").ColumnWidth = 8
$unit.Columns("C:C").ColumnWidth = 14
$unit.Columns("D:D").ColumnWidth = 21
$unit.Columns("E:E").ColumnWidth = 32

$unit.Activate()
$unit.Range("E10").Select()

# ---------------------------------------------------------------------
# 6. Final active sheet/tab: eia_clean (matches activeTab=1 / tabSelected)
# ---------------------------------------------------------------------
$eia = $wb.Worksheets.Item("eia_clean")
$eia.Activate()
$eia.Range("C1").Select()
